# Update the "Global" data sheet to reflect the new Fiori Launchpad URL
# for the environment, and adjust the view accordingly.

$wb = $excel.ActiveWorkbook
$wsGlobal = $wb.Worksheets.Item("Global")
$wsOverview = $wb.Worksheets.Item("GLOverview")

# Update the URL value held in cell B2
$wsGlobal.Range("B2").Value = "https://sap-hana-vpn.mfdemoportal.com:44300/sap/bc/ui2/flp?sap-client=100&sap-language=EN#Shell-home"

# Resize column B so the longer URL text fits (auto-fit based on content)
$wsGlobal.Activate()
$wsGlobal.Columns("B:B").ColumnWidth = 89.17

# Move the active selection to B1 on the Global sheet
$wsGlobal.Range("B1").Select()

# Restore the originally active sheet (GLOverview)
$wsOverview.Activate()
